# Resort the sheet tabs: move "总计" (the summary sheet) so that it
# becomes the first sheet, ahead of "2021-Q1" (the detail sheet).
$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("总计")

# Move "总计" to be before the current first sheet, i.e. make it sheet #1.
$wsSummary.Move($wb.Worksheets.Item(1))
